$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 00:12"

# --- Reorder country rows (names swap; some rows carry their stats with them) ---
# Fiyi (row202) / Dominica (row203) swap names only - stats identical, so safe either way
$ws.Cells.Item(202, 1).Value = "Dominica"
$ws.Cells.Item(203, 1).Value = "Fiyi"

# Groenlandia (row207) / Islas Malvinas (row208) swap names only - stats identical
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"
$ws.Cells.Item(208, 1).Value = "Groenlandia"

# Papua Nueva Guinea (row213) / Islas Virgenes Britanicas (row214) swap fully (names + stats)
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0

# --- Updated COVID stats (refreshed data pull) ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 2353775
$ws.Cells.Item(4, 3).Value = 23197
$ws.Cells.Item(4, 4).Value = 977683
$ws.Cells.Item(4, 5).Value = 1253854
$ws.Cells.Item(4, 7).Value = 258
$ws.Cells.Item(4, 8).Value = 122238

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 1084883
$ws.Cells.Item(5, 3).Value = 14744
$ws.Cells.Item(5, 5).Value = 491089
$ws.Cells.Item(5, 7).Value = 550
$ws.Cells.Item(5, 8).Value = 50608

# Row 50: Barein
$ws.Cells.Item(50, 2).Value = 21764
$ws.Cells.Item(50, 3).Value = 433
$ws.Cells.Item(50, 4).Value = 16419
$ws.Cells.Item(50, 5).Value = 5282

# Row 90: Bulgaria
$ws.Cells.Item(90, 2).Value = 3905
$ws.Cells.Item(90, 3).Value = 33
$ws.Cells.Item(90, 4).Value = 2074
$ws.Cells.Item(90, 5).Value = 1632

# Row 160: Surinam
$ws.Cells.Item(160, 2).Value = 314
$ws.Cells.Item(160, 3).Value = 12
$ws.Cells.Item(160, 4).Value = 106
$ws.Cells.Item(160, 5).Value = 200
